$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (29) below the last existing row (28), carrying the
# same centered formatting used by the rest of the data rows.
$ws.Range("A28:C28").Copy()
$ws.Range("A29:C29").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Column A holds a date-like label ("2025/12/08") that must stay plain text
# (matching the existing "2025/MM/DD" entries above it) instead of being
# auto-converted to a date serial number, so the cell is switched to Text
# format before the value is entered.
$ws.Range("A29").NumberFormat = "@"
$ws.Range("A29").Value = "2025/12/08"

$ws.Range("B29").Value = "逃离鸭科夫"
$ws.Range("C29").Value = 1349

# Re-apply the shared centered/General formatting (copied from an existing
# data cell) on top of A29 so it ends up with the exact same style as the
# rest of column A, without altering the text value just entered.
$ws.Range("B28").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
